# Apply latest crypto price/volume snapshot to Sheet1 (columns D, E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.184.15'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '3.610.90'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'579.96"
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').Value = "'175.32"
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('D7').Value = '3.602.94'
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('D8').Value = "'0.608"
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('D10').Value = "'0.196"
$ws.Range('E10').Value = '  -5.25%  '
$ws.Range('D11').Value = "'6.78"
$ws.Range('E11').Value = '  +22.83%  '
$ws.Range('D12').Value = "'0.602"
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = "'48.22"
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '4.188.07'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = "'667.54"
$ws.Range('E16').Value = '  -3.86%  '
$ws.Range('D17').Value = "'8.85"
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '3.606.12'
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').Value = '70.169.84'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = "'17.69"
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').Value = "'11.34"
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').Value = "'0.928"
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = "'17.03"
$ws.Range('E24').Value = '  -2.77%  '
$ws.Range('D25').Value = "'99.47"
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('D26').Value = "'3.89"
$ws.Range('E26').Value = '  -2.89%  '
$ws.Range('D27').Value = "'2.77"
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = "'9.90"
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = "'34.45"
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('E33').Value = '  -6.71%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -4.16%  '
$ws.Range('D36').Value = "'574.42"
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('D37').Value = "'11.01"
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D39').Value = "'57.88"
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = '3.572.65'
$ws.Range('E41').Value = '  -2.60%  '
$ws.Range('D42').Value = "'0.0453"
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('D43').Value = "'0.139"
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('D45').Value = "'34.45"
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('E46').Value = '  -3.89%  '
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('D48').Value = "'2.87"
$ws.Range('E48').Value = '  +4.73%  '
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('D50').Value = "'136.14"
$ws.Range('E50').Value = '  +1.74%  '
$ws.Range('D51').Value = "'2.93"
$ws.Range('E51').Value = '  +2.28%  '
